# Trade #36 closed at 2026-02-17 08:32:51 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" roll-up figures to reflect the
# newly closed trade, and appends the new trade row (#36 / row 37) to both
# the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.57    # Current Capital
$summary.Range("B4").Value = -0.43      # Total P&L $
$summary.Range("B5").Value = -0.24      # Total P&L %
$summary.Range("B6").Value = 36         # Total Trades
$summary.Range("B7").Value = 11         # Winning Trades
$summary.Range("B9").Value = 30.56      # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.57       # Capital
$status.Range("D4").Value = 36          # Trades
$status.Range("E4").Value = -0.43       # P&L $
$status.Range("F4").Value = -0.43       # P&L %
$status.Range("G4").Value = 30.56       # Win Rate %

# ---------------------------------------------------------------------------
# Append the new trade row (row 37) to a sheet
# ---------------------------------------------------------------------------
function Add-Trade37 {
    param($ws)

    $row = 37

    $ws.Cells.Item($row, 1).Value = 36

    # Date-looking text must stay as plain text, not get auto-converted to a
    # date serial by Excel's smart entry. Force text format, assign, then
    # reset the style back to Normal so no extra style gets attached.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:32:44"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.61
    $ws.Cells.Item($row, 7).Value = 0.62
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 1.6393
    $ws.Cells.Item($row, 10).Value = 0.01
    $ws.Cells.Item($row, 11).Value = 99.57
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Add-Trade37 $wb.Worksheets.Item("All Trades")
Add-Trade37 $wb.Worksheets.Item("MarketMaking")
